$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 379
$lastHyperlinkRow = 43
$hyperlinkCols = @(19, 20, 21, 22, 23, 24, 25)  # S, T, U, V, W, X, Y

for ($r = 2; $r -le $lastRow; $r++) {
    # Update "Förändrad" (changed) date in column C
    $ws.Cells.Item($r, 3).Value2 = 45186

    if ($r -le $lastHyperlinkRow) {
        $beteckning = $ws.Cells.Item($r, 1).Value2

        foreach ($c in $hyperlinkCols) {
            $cell = $ws.Cells.Item($r, $c)
            $f = $cell.Formula
            if ($f -ne "") {
                $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $beteckning + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
